# Applies corrected Diebold-Mariano p-values and statistics to the
# P_valores and Estadisticos_DM sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "P_valores" ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.1526727808725477
$wsP.Range("D2").Value = 0.3974693935216091
$wsP.Range("E2").Value = 0.3847456397084215
$wsP.Range("F2").Value = 0.1459621063166057

$wsP.Range("B3").Value = 0.1526727808725477
$wsP.Range("D3").Value = 0.4589534655230383
$wsP.Range("E3").Value = 0.2941538800409709
$wsP.Range("F3").Value = 0.5610034795744343

$wsP.Range("B4").Value = 0.3974693935216091
$wsP.Range("C4").Value = 0.4589534655230383
$wsP.Range("E4").Value = 0.6123118183722083
$wsP.Range("F4").Value = 0.7701633429726855

$wsP.Range("B5").Value = 0.3847456397084215
$wsP.Range("C5").Value = 0.2941538800409709
$wsP.Range("D5").Value = 0.6123118183722083
$wsP.Range("F5").Value = 0.3706408452925514

$wsP.Range("B6").Value = 0.1459621063166057
$wsP.Range("C6").Value = 0.5610034795744343
$wsP.Range("D6").Value = 0.7701633429726855
$wsP.Range("E6").Value = 0.3706408452925514

# --- Sheet "Estadisticos_DM" ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -1.512406264111316
$wsE.Range("D2").Value = -0.8728489151157318
$wsE.Range("E2").Value = -0.8972695369579118
$wsE.Range("F2").Value = -1.539552960120117

$wsE.Range("B3").Value = 1.512406264111316
$wsE.Range("D3").Value = 0.7615714069895181
$wsE.Range("E3").Value = 1.089918960120735
$wsE.Range("F3").Value = 0.5955148722530386

$wsE.Range("B4").Value = 0.8728489151157318
$wsE.Range("C4").Value = -0.7615714069895181
$wsE.Range("E4").Value = 0.5183456872182491
$wsE.Range("F4").Value = -0.2978883076503569

$wsE.Range("B5").Value = 0.8972695369579118
$wsE.Range("C5").Value = -1.089918960120735
$wsE.Range("D5").Value = -0.5183456872182491
$wsE.Range("F5").Value = -0.9249905528605229

$wsE.Range("B6").Value = 1.539552960120117
$wsE.Range("C6").Value = -0.5955148722530386
$wsE.Range("D6").Value = 0.2978883076503569
$wsE.Range("E6").Value = 0.9249905528605229

$wb.Save()
